$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1249
$ws1.Range("F3").Value = 657
$ws1.Range("F4").Value = 358
$ws1.Range("F5").Value = 5082
$ws1.Range("F6").Value = 539
$ws1.Range("F7").Value = 9703
$ws1.Range("F9").Value = 541
$ws1.Range("F10").Value = 96
$ws1.Range("F11").Value = 27
$ws1.Range("F12").Value = 716
$ws1.Range("F13").Value = 79

# Sheet 2: 演出 (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 14
$ws2.Range("F4").Value = 12
$ws2.Range("F5").Value = 8
$ws2.Range("F6").Value = 3

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1249
$ws4.Range("F3").Value = 657
$ws4.Range("F4").Value = 358
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 5082
$ws4.Range("F8").Value = 539
$ws4.Range("F9").Value = 12
$ws4.Range("F10").Value = 9703
$ws4.Range("F11").Value = 252
$ws4.Range("F12").Value = 541
$ws4.Range("F13").Value = 96
$ws4.Range("F14").Value = 8
$ws4.Range("F16").Value = 27
$ws4.Range("F17").Value = 716
$ws4.Range("F19").Value = 79

$wb.Save()
